# Update cryptocurrency price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '28.460.14'
$ws.Cells.Item(2, 5).Value = '  +0.33%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.828.04'
$ws.Cells.Item(3, 5).Value = '  +0.13%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.23%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '315.24'
$ws.Cells.Item(5, 5).Value = '  -0.84%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.27%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5150'
$ws.Cells.Item(7, 5).Value = '  -3.49%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3930'
$ws.Cells.Item(8, 5).Value = '  -1.55%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07661'
$ws.Cells.Item(9, 5).Value = '  +1.57%  '

# Row 10
$ws.Cells.Item(10, 2).Value = 'Polygon'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.111'
$ws.Cells.Item(10, 5).Value = '  +0.71%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'OKB'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '41.70'
$ws.Cells.Item(11, 5).Value = '  -0.29%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '21.08'
$ws.Cells.Item(12, 5).Value = '  +2.23%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.302'
$ws.Cells.Item(13, 5).Value = '  -0.06%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.001'
$ws.Cells.Item(14, 5).Value = '  +0.17%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.544'
$ws.Cells.Item(15, 5).Value = '  -1.08%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.821.64'
$ws.Cells.Item(16, 5).Value = '  +0.03%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '93.68'
$ws.Cells.Item(17, 5).Value = '  +4.54%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001104'
$ws.Cells.Item(18, 5).Value = '  +3.31%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06684'
$ws.Cells.Item(19, 5).Value = '  +1.59%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.68'
$ws.Cells.Item(20, 5).Value = '  +1.27%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.40%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.149'
$ws.Cells.Item(22, 5).Value = '  +2.08%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '28.483.76'
$ws.Cells.Item(23, 5).Value = '  +0.32%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.17'
$ws.Cells.Item(24, 5).Value = '  -0.15%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.257'

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '20.79'
$ws.Cells.Item(26, 5).Value = '  +1.42%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '156.93'
$ws.Cells.Item(27, 5).Value = '  +0.08%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.035.37'
$ws.Cells.Item(28, 5).Value = '  +0.38%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.405'
$ws.Cells.Item(29, 5).Value = '  +0.66%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '124.70'
$ws.Cells.Item(30, 5).Value = '  +1.13%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.115'
$ws.Cells.Item(31, 5).Value = '  +0.50%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.1084'
$ws.Cells.Item(32, 5).Value = '  -1.32%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.665'
$ws.Cells.Item(33, 5).Value = '  +1.33%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.667'
$ws.Cells.Item(34, 5).Value = '  -0.37%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.07012'
$ws.Cells.Item(35, 5).Value = '  -4.08%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.2209'
$ws.Cells.Item(36, 5).Value = '  -1.21%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '8.923'
$ws.Cells.Item(37, 5).Value = '  +2.66%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02322'
$ws.Cells.Item(38, 5).Value = '  +0.41%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '5.167'
$ws.Cells.Item(39, 5).Value = '  -1.04%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.6267'
$ws.Cells.Item(40, 5).Value = '  +0.69%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.71%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.176'
$ws.Cells.Item(42, 5).Value = '  -1.43%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.000'
$ws.Cells.Item(43, 5).Value = '  +0.23%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.389'
$ws.Cells.Item(44, 5).Value = '  -1.36%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.46'
$ws.Cells.Item(45, 5).Value = '  -0.07%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5899'

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.709'
$ws.Cells.Item(47, 5).Value = '  +0.26%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '124.67'
$ws.Cells.Item(48, 5).Value = '  -0.33%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.979'
$ws.Cells.Item(49, 5).Value = '  +1.25%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.197'
$ws.Cells.Item(50, 5).Value = '  +0.75%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06927'
$ws.Cells.Item(51, 5).Value = '  +0.57%  '

